$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 13. This pushes the existing
# "Programa resumido:" ... "Requisitos:" label block (and the LOQ4098 row)
# down by two rows, making room for a dedicated two-row "Docentes
# responsaveis" block (one row per teacher) right under row 12.
$ws.Rows.Item(13).Resize(2).Insert()

# Row 10 (Objetivos:) - replace the misplaced teacher name with the real
# objectives text.
$ws.Range("B10:C10").Value = 'Dotar os alunos de conhecimentos básicos sobre o uso da corrente elétrica como agente de redução de metais com aplicação nos processos eletrolíticos de extração, refino, revestimento e formação de peças.'

# New rows 13/14 under "Docentes responsáveis:" (row 12) - one teacher per row.
$ws.Range("B13:C13").Value = '144651 - Antonio Fernando Sartori'
$ws.Range("B14:C14").Value = '3577649 - Carlos Angelo Nunes'

# Row 15 (Programa resumido:) - replace stray date with the real short syllabus.
$ws.Range("B15:C15").Value = 'Fundamentos da Eletrólise Aplicados a Eletrodeposição de Metais. Aplicações Industriais da Eletrodeposição.'

# Row 17 (Programa:) - replace stray teacher name with the real syllabus text.
$ws.Range("B17:C17").Value = 'Aplicações da Eletrodeposição. Fundamentos da Eletrólise Aplicada a Eletrodeposição de Metais: Potencial de Eletrodo. Dupla Camada Elétrica. Lei de Nernst. Células Galvânicas e Eletrolíticas. Processos Catódicos e Anódicos: Cinética Eletroquímica e Polarização dos Eletrodos. Corrente de Troca. Processos Galvanotécnicos e Eletrometalurgicos: Eletrorrevestimento, Eletroformação de Peças, Eletrorrefino, Eletroextração e Eletrodeposição de ligas. Aplicações industriais: Eletrodeposição de níquel, cromo, cobre, ouro, prata, estanho e nióbio. Eletroextração de Alumínio.'

# Row 20 (Método:) - replace stray teacher name with the grading method text.
$ws.Range("B20:C20").Value = 'Serão usadas duas notas P1 e P2. A P1 será uma prova escrita e a P2 será a soma de uma nota de seminário e uma nota de trabalho escrito.'

# Row 21 (Critério:) - replace the method text (now duplicated via the shift)
# with the actual grading criteria text.
$ws.Range("B21:C21").Value = 'Média Final(MF) = (P1 + 2P2)/3 MF menor que 3,0 - reprovado. MF maior,ou igual, a 3,0 até menor que 5,0 - recuperação. MF maior, ou igual, a 5,0 - aprovado.'

# Row 22 (Norma de recuperação:) - replace the criteria text with the real
# recovery-exam rule text.
$ws.Range("B22:C22").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'

# Row 23 (Bibliografia:) - replace the recovery-rule text with the real
# bibliography list.
$ws.Range("B23:C23").Value = '1. Modern Electrochemistry, J.O.Bockris, A. K. N. Reddy, M. G. Aldeco, Plenum Publishing Corporation, 2a ed, 2000.2. Comprehensive Treatise of Eletrochemistry, Bockris, Kluwer Academic Pub, 1981.3. Modern Electroplating, F. Lowenhein, John Wiley and Sons, 3a ed, 1974.4. Modern electroplating, M. Schlesinger, M. Paunovic, Wiley-Interscience, 4a ed, 2000. 5. Fundamentals of Electochemical Deposition (Electrochemical Society Series), M. Paunovic, M. Schlesinger, Wiley-Interscience, 1998. 6.  ASM Handbook: Coating and Coating Processes for metals (Materials data series), J. Lindsay, ASM International, 1998.  7. Eletrodeposition of alloy, A.Brenner, Academic Press, 1963. 8. Denaro, A.R. Fundamentos de Eletroquímica. Editora Edgard Blucher, 1974.9. Electroplating, N. Kanani, Elsevier, 2004.10. Electrodeposition, J. W. Dini, Noyes Publications, 1993.'

# Rows 13/14 are brand-new plain data rows (no bold label in column A), and
# rows 15-23 must keep the same row heights they had before the insert
# (AutoFit clears the custom height that Insert() may have copied down).
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Rows.Item(14).EntireRow.AutoFit()
